# The edit described by the diff is a pure re-ordering of the data rows
# (rows 2-13 on the "Artfynd" sheet) — every row's full set of columns
# (A:AY) moves intact to a new row position, keyed by the "Id" value in
# column A. No cell values themselves change, only which row they sit on.
#
# Mapping of new row -> old row (by matching column A / "Id"):
#   2<-11  3<-7  4<-6  5<-9  6<-5  7<-8  8<-12  9<-4  10<-13  11<-3  12<-2  13<-10

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstDataRow = 2
$lastDataRow  = 13

$usedRange  = $ws.UsedRange
$lastCol    = $usedRange.Columns.Count
$lastColRef = $ws.Cells.Item(1, $lastCol).Address($false, $false)
$lastColLetter = ($lastColRef -replace '[0-9]', '')

# new row -> old row
$mapping = [ordered]@{
    2  = 11
    3  = 7
    4  = 6
    5  = 9
    6  = 5
    7  = 8
    8  = 12
    9  = 4
    10 = 13
    11 = 3
    12 = 2
    13 = 10
}

# Use a scratch area, well below the real data, to stage a verbatim copy of
# every source row before anything gets overwritten (this is a permutation,
# so a row can be both a source and a destination).
$stageOffset = 1000

for ($r = $firstDataRow; $r -le $lastDataRow; $r++) {
    $srcRange   = $ws.Range("A$r`:$lastColLetter$r")
    $stageRow   = $r + $stageOffset
    $stageRange = $ws.Range("A$stageRow`:$lastColLetter$stageRow")
    $srcRange.Copy($stageRange)
}

# Clear the destination rows first: Copy-paste here only overwrites cells
# that actually hold a value in the source, so stale content left behind
# from cells that are empty in the source (but populated in the old
# destination) would otherwise survive.
for ($r = $firstDataRow; $r -le $lastDataRow; $r++) {
    $ws.Range("A$r`:$lastColLetter$r").ClearContents()
}

foreach ($newRow in $mapping.Keys) {
    $oldRow     = $mapping[$newRow]
    $stageRow   = $oldRow + $stageOffset
    $stageRange = $ws.Range("A$stageRow`:$lastColLetter$stageRow")
    $destRange  = $ws.Range("A$newRow`:$lastColLetter$newRow")
    $stageRange.Copy($destRange)
}

# Clean up the scratch area so the sheet dimensions go back to normal.
for ($r = $firstDataRow; $r -le $lastDataRow; $r++) {
    $stageRow = $r + $stageOffset
    $ws.Range("A$stageRow`:$lastColLetter$stageRow").ClearContents()
}
